$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Add the "comodin" product row details: Responsable = Agustina, Avance = 100%
$ws.Range("B36").Value = "Agustina"
$ws.Range("C36").Value = 1
$ws.Range("C36").NumberFormat = "0%"

# Update the selection saved in the sheet view
$ws.Range("B41").Select()
